$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three data rows (2007年, 2008年, 2009年) — everything
# below shifts up to fill the gap.
$ws.Rows("2:4").Delete()

# Append a new trailing row for 2021年 (only the aggregate column J is
# populated, matching the pattern already used for 2019年/2020年).
$ws.Range("A13").Value = "2021年"
$ws.Range("J13").Value = 2050

# Give the new year label cell (A13) the same style as the rest of
# column A (centered/bordered header-ish style used for year labels).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
